$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "stimulus_duration" column (old column E); this shifts F:J left to E:I
$ws.Range("E1").EntireColumn.Delete()

# Update current_font values from font file paths to font display names
$ws.Range("A2").Value = "Georgia"
$ws.Range("A3").Value = "RobotoFlex"
$ws.Range("A4").Value = "Neue Frutiger World"
$ws.Range("A5").Value = "Georgia"
$ws.Range("A6").Value = "RobotoFlex"
$ws.Range("A7").Value = "Neue Frutiger World"

# Update word_trial_count / nonword_trial_count for the lexical training row (staircase change)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1

# Update the selected cell in the sheet view
[void]$ws.Range("A5").Select()
